# Generate Report for Handback
#
# Updates the "6dfee820-d28c-4894-8c6a-17eaf39aee9a.md" handback row (row 7)
# on both locale sheets (zh-cn, de-de): the handback was reprocessed and
# found to be against a stale source version, so the "Latest Target File",
# "Latest Handback File", "Latest Handback DateTime" and "Error Detail"
# columns now get populated (they were placeholders before).

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/56256cc2e08f4d1d74f355a506789f7e0cfc05c3/e2e/6dfee820-d28c-4894-8c6a-17eaf39aee9a.md"
$handbackDisplay = "6dfee820-d28c-4894-8c6a-17eaf39aee9a.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/368d90ed644a1e695a6b797774157a4417ec136e/e2e/6dfee820-d28c-4894-8c6a-17eaf39aee9a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/56256cc2e08f4d1d74f355a506789f7e0cfc05c3/e2e/6dfee820-d28c-4894-8c6a-17eaf39aee9a.md."

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = $handbackDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $handbackUrl, "", "", $handbackDisplay)

$wsZh.Range("J7").Value = "6dfee820-d28c-4894-8c6a-17eaf39aee9a.9ff7fb3f4a1d7a00fa92ff462b216592396a40f6.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-21 23:03:20"
$wsZh.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = $handbackDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $handbackUrl, "", "", $handbackDisplay)

$wsDe.Range("J7").Value = "6dfee820-d28c-4894-8c6a-17eaf39aee9a.9ff7fb3f4a1d7a00fa92ff462b216592396a40f6.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-21 23:03:26"
$wsDe.Range("P7").Value = $errorDetail
